# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Update the "time_taken" timestamps on the existing "data" sheet
# -----------------------------------------------------------------
$data = $wb.Worksheets.Item("data")

$data.Range("F2").Value = "2021-10-05 14:35:35.978160"
$data.Range("F3").Value = "2021-10-05 14:35:35.978167"
$data.Range("F4").Value = "2021-10-05 14:35:35.978171"
$data.Range("F5").Value = "2021-10-05 14:35:35.978174"
$data.Range("F6").Value = "2021-10-05 14:35:35.978177"
$data.Range("F7").Value = "2021-10-05 14:35:35.978179"
$data.Range("F8").Value = "2021-10-05 14:35:35.978182"

# -----------------------------------------------------------------
# 2) Add a new "metadata" worksheet right after "data"
# -----------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (B1:G1) - reuse the bold/centered/bordered header style
# already used by the "data" tab so no duplicate style gets created.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Index cell A2, styled like the "data" tab's A-column index cells
$meta.Range("A2").Value = 0
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Data row
$meta.Range("B2").Value = "Renal Tubulointerstitial Disease"
$meta.Range("C2").Value = 199

$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.0"

$meta.Range("E2").Value = "2021-01-16T11:32:11.361557Z"
$meta.Range("F2").Value = "2021-10-05 14:35:35.974525"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/199/?format=json"

$excel.CutCopyMode = $false
$data.Activate()
[void]$data.Range("A1").Select()
